$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing data rows (old 2..22) down to (3..23)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new data point
$ws.Range("A2").Value = 0.0265726372599601
$ws.Range("B2").Value = -0.0716239511966705
$ws.Range("C2").Value = -0.0607810914516449

# After the insert, old row 21 and old row 22 live at rows 22 and 23.
# The target sheet only keeps data through row 21, so drop the trailing two rows.
$ws.Range("A22:A23").EntireRow.Delete()
